$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 16:50"

# 2) Move "Republica Dominicana" up to right after "Panama" (row 47),
#    pushing Mexico / Singapur / Argentina / Serbia / Eslovenia down one row each
#    (rows 47-52). Each row is rewritten in full with the country name and
#    figures it ends up with after the move.
$ws.Cells.Item(47,1).Value = "Republica Dominicana"
$ws.Cells.Item(47,2).Value = 859
$ws.Cells.Item(47,3).Value = 140
$ws.Cells.Item(47,4).Value = 3
$ws.Cells.Item(47,5).Value = 817
$ws.Cells.Item(47,6).Value = 0
$ws.Cells.Item(47,7).Value = 11
$ws.Cells.Item(47,8).Value = 39

$ws.Cells.Item(48,1).Value = "Mexico"
$ws.Cells.Item(48,2).Value = 848
$ws.Cells.Item(48,3).Value = 131
$ws.Cells.Item(48,4).Value = 4
$ws.Cells.Item(48,5).Value = 828
$ws.Cells.Item(48,6).Value = 1
$ws.Cells.Item(48,7).Value = 4
$ws.Cells.Item(48,8).Value = 16

$ws.Cells.Item(49,1).Value = "Singapur"
$ws.Cells.Item(49,2).Value = 844
$ws.Cells.Item(49,3).Value = 42
$ws.Cells.Item(49,4).Value = 212
$ws.Cells.Item(49,5).Value = 629
$ws.Cells.Item(49,6).Value = 19
$ws.Cells.Item(49,7).Value = 1
$ws.Cells.Item(49,8).Value = 3

$ws.Cells.Item(50,1).Value = "Argentina"
$ws.Cells.Item(50,2).Value = 745
$ws.Cells.Item(50,3).Value = 0
$ws.Cells.Item(50,4).Value = 72
$ws.Cells.Item(50,5).Value = 654
$ws.Cells.Item(50,6).Value = 0
$ws.Cells.Item(50,7).Value = 0
$ws.Cells.Item(50,8).Value = 19

$ws.Cells.Item(51,1).Value = "Serbia"
$ws.Cells.Item(51,2).Value = 741
$ws.Cells.Item(51,3).Value = 82
$ws.Cells.Item(51,4).Value = 42
$ws.Cells.Item(51,5).Value = 686
$ws.Cells.Item(51,6).Value = 25
$ws.Cells.Item(51,7).Value = 3
$ws.Cells.Item(51,8).Value = 13

$ws.Cells.Item(52,1).Value = "Eslovenia"
$ws.Cells.Item(52,2).Value = 730
$ws.Cells.Item(52,3).Value = 46
$ws.Cells.Item(52,4).Value = 10
$ws.Cells.Item(52,5).Value = 709
$ws.Cells.Item(52,6).Value = 23
$ws.Cells.Item(52,7).Value = 2
$ws.Cells.Item(52,8).Value = 11

# 3) Standalone daily numeric refreshes for other countries
# Estados Unidos (row 4)
$ws.Cells.Item(4,2).Value = 123898
$ws.Cells.Item(4,3).Value = 320
$ws.Cells.Item(4,5).Value = 118429

# Reino Unido (row 11)
$ws.Cells.Item(11,5).Value = 18159
$ws.Cells.Item(11,7).Value = 209
$ws.Cells.Item(11,8).Value = 1228

# Noruega (row 20)
$ws.Cells.Item(20,2).Value = 4239
$ws.Cells.Item(20,3).Value = 224
$ws.Cells.Item(20,5).Value = 4207

# Israel (row 23)
$ws.Cells.Item(23,5).Value = 3761
$ws.Cells.Item(23,7).Value = 3
$ws.Cells.Item(23,8).Value = 15

# Chequia (row 25)
$ws.Cells.Item(25,2).Value = 2716
$ws.Cells.Item(25,3).Value = 85
$ws.Cells.Item(25,5).Value = 2692

# Pakistan (row 36)
$ws.Cells.Item(36,5).Value = 1483
$ws.Cells.Item(36,7).Value = 2
$ws.Cells.Item(36,8).Value = 14

# Bulgaria (row 73)
$ws.Cells.Item(73,2).Value = 346
$ws.Cells.Item(73,3).Value = 15
$ws.Cells.Item(73,4).Value = 14
$ws.Cells.Item(73,5).Value = 324
$ws.Cells.Item(73,6).Value = 9

# Sri Lanka (row 104)
$ws.Cells.Item(104,4).Value = 11
$ws.Cells.Item(104,5).Value = 103

# Zambia (row 138)
$ws.Cells.Item(138,2).Value = 29
$ws.Cells.Item(138,3).Value = 1
$ws.Cells.Item(138,5).Value = 29
